$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.334.67"
$ws.Range("E2").Value = "  +4.19%  "

$ws.Range("D3").Value = "3.485.69"
$ws.Range("E3").Value = "  +5.30%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "'555.65"
$ws.Range("E5").Value = "  +7.19%  "

$ws.Range("D6").Value = "'184.34"
$ws.Range("E6").Value = "  +7.62%  "

$ws.Range("D7").Value = "'0.637"
$ws.Range("E7").Value = "  +8.42%  "

$ws.Range("D8").Value = "3.482.28"
$ws.Range("E8").Value = "  +5.23%  "

$ws.Range("E9").Value = "  +0.05%  "

$ws.Range("D10").Value = "'0.632"
$ws.Range("E10").Value = "  +5.29%  "

$ws.Range("D11").Value = "'0.152"
$ws.Range("E11").Value = "  +15.11%  "

$ws.Range("D12").Value = "'54.16"
$ws.Range("E12").Value = "  +2.84%  "

$ws.Range("D13").Value = "'0.0000270"
$ws.Range("E13").Value = "  +6.35%  "

$ws.Range("D14").Value = "'9.26"
$ws.Range("E14").Value = "  +3.66%  "

$ws.Range("D15").Value = "4.053.58"
$ws.Range("E15").Value = "  +5.48%  "

$ws.Range("D16").Value = "3.490.90"
$ws.Range("E16").Value = "  +5.41%  "

$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "'0.121"
$ws.Range("E17").Value = "  +4.32%  "

$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").Value = "'18.41"
$ws.Range("E18").Value = "  +5.95%  "

$ws.Range("D19").Value = "66.340.79"
$ws.Range("E19").Value = "  +4.59%  "

$ws.Range("D20").Value = "'11.99"
$ws.Range("E20").Value = "  +8.20%  "

$ws.Range("D21").Value = "'0.990"
$ws.Range("E21").Value = "  +4.37%  "

$ws.Range("D22").Value = "'421.05"
$ws.Range("E22").Value = "  +12.97%  "

$ws.Range("D23").Value = "'4.04"
$ws.Range("E23").Value = "  +11.30%  "

$ws.Range("D24").Value = "'85.91"
$ws.Range("E24").Value = "  +6.02%  "

$ws.Range("E25").Value = "  -0.92%  "

$ws.Range("D26").Value = "'10.86"
$ws.Range("E26").Value = "  -3.85%  "

$ws.Range("D27").Value = "'2.89"
$ws.Range("E27").Value = "  +8.38%  "

$ws.Range("D28").Value = "'12.22"
$ws.Range("E28").Value = "  +9.70%  "

$ws.Range("D29").Value = "'6.07"
$ws.Range("E29").Value = "  -1.57%  "

$ws.Range("D30").Value = "'9.08"
$ws.Range("E30").Value = "  +12.15%  "

$ws.Range("D31").Value = "'30.10"
$ws.Range("E31").Value = "  +5.47%  "

$ws.Range("D32").Value = "'628.49"
$ws.Range("E32").Value = "  +1.41%  "

$ws.Range("D33").Value = "'6.57"
$ws.Range("E33").Value = "  +3.44%  "

$ws.Range("D34").Value = "'11.70"
$ws.Range("E34").Value = "  +5.43%  "

$ws.Range("E35").Value = "  +5.78%  "

$ws.Range("D36").Value = "'59.90"
$ws.Range("E36").Value = "  +3.78%  "

$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "'0.147"
$ws.Range("E37").Value = "  +20.15%  "

$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0807"
$ws.Range("E38").Value = "  +11.97%  "

$ws.Range("E39").Value = "  -0.03%  "

$ws.Range("D40").Value = "'37.51"
$ws.Range("E40").Value = "  +5.59%  "

$ws.Range("D41").Value = "'0.384"
$ws.Range("E41").Value = "  +2.68%  "

$ws.Range("D42").Value = "'3.50"
$ws.Range("E42").Value = "  +16.37%  "

$ws.Range("D43").Value = "3.114.80"
$ws.Range("E43").Value = "  +8.11%  "

$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  +0.20%  "

$ws.Range("E45").Value = "  -0.68%  "

$ws.Range("D46").Value = "'2.83"
$ws.Range("E46").Value = "  +10.39%  "

$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0413"
$ws.Range("E47").Value = "  +5.55%  "

$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "'3.24"
$ws.Range("E48").Value = "  +9.13%  "

$ws.Range("E49").Value = "  +3.05%  "

$ws.Range("D50").Value = "'0.133"
$ws.Range("E50").Value = "  +7.13%  "

$ws.Range("D51").Value = "'140.61"
$ws.Range("E51").Value = "  +3.43%  "
